$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.983.63"
$ws.Range("E2").Value = "  +0.55%  "
$ws.Range("D3").Value = "1.739.59"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.44"
$ws.Range("E5").Value = "  +4.51%  "
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5026"
$ws.Range("E7").Value = "  -1.80%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2730"
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06191"
$ws.Range("E9").Value = "  +1.39%  "
$ws.Range("B10").Value = "WrappedEther"
$ws.Range("C10").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D10").Value = "1.744.16"
$ws.Range("E10").Value = "  +0.33%  "
$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07259"
$ws.Range("E11").Value = "  +1.45%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.6538"
$ws.Range("E12").Value = "  +2.85%  "
$ws.Range("E13").Value = "  +1.60%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.738"
$ws.Range("E14").Value = "  +3.35%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.73"
$ws.Range("E15").Value = "  +0.84%  "
$ws.Range("E16").Value = "  +0.05%  "
$ws.Range("E17").Value = "  +0.13%  "
$ws.Range("D18").Value = "25.998.62"
$ws.Range("E18").Value = "  +0.57%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.89"
$ws.Range("E19").Value = "  +1.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000006848"
$ws.Range("E20").Value = "  +2.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.606"
$ws.Range("E21").Value = "  +8.54%  "
$ws.Range("D22").Value = "1.963.11"
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.755"
$ws.Range("E23").Value = "  +1.33%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.403"
$ws.Range("E24").Value = "  +3.58%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "134.10"
$ws.Range("E25").Value = "  -3.42%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.498"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.26"
$ws.Range("E27").Value = "  +1.03%  "
$ws.Range("E28").Value = "  +2.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "105.27"
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.991"
$ws.Range("E30").Value = "  +1.21%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08148"
$ws.Range("E31").Value = "  -2.29%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.703"
$ws.Range("E32").Value = "  +1.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04735"
$ws.Range("E33").Value = "  +4.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.666"
$ws.Range("E34").Value = "  +0.22%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9962"
$ws.Range("E35").Value = "  +1.70%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6140"
$ws.Range("E36").Value = "  -0.56%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.751"
$ws.Range("E37").Value = "  +2.40%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01619"
$ws.Range("E38").Value = "  +1.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.944"
$ws.Range("E39").Value = "  +1.69%  "
$ws.Range("E40").Value = "  +0.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "100.90"
$ws.Range("E41").Value = "  +3.36%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8035"
$ws.Range("E42").Value = "  +9.53%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.3914"
$ws.Range("E43").Value = "  +2.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.013"
$ws.Range("E44").Value = "  +1.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1173"
$ws.Range("E45").Value = "  +4.40%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.372"
$ws.Range("E46").Value = "  +3.46%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.79"
$ws.Range("E47").Value = "  +2.05%  "
$ws.Range("E48").Value = "  +0.61%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.91"
$ws.Range("E49").Value = "  +1.59%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.663"
$ws.Range("E50").Value = "  +1.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3482"
$ws.Range("E51").Value = "  +2.14%  "
